# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new F-value for sheet "展览"
$updates1 = @{
    5  = 90
    7  = 1258
    8  = 1533
    9  = 340
    10 = 388
    12 = 147
    13 = 163
    15 = 106
    18 = 322
    19 = 1730
    26 = 4169
    28 = 268
    29 = 1085
    32 = 529
    34 = 243
}

# Row -> new F-value for sheet "全部类型"
$updates4 = @{
    5  = 90
    7  = 1258
    8  = 1533
    9  = 340
    10 = 388
    12 = 147
    13 = 163
    15 = 106
    18 = 322
    19 = 1730
    26 = 4170
    28 = 268
    29 = 1085
    32 = 529
    34 = 243
}

foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
